$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 0.5603867769241333
$ws.Cells.Item(3, 4).Value = 0.3992350995540619
$ws.Cells.Item(4, 4).Value = 0.9112033843994141
$ws.Cells.Item(5, 4).Value = 0.6715229749679565
$ws.Cells.Item(6, 4).Value = 0.198024719953537
$ws.Cells.Item(7, 4).Value = 0.9465768337249756
$ws.Cells.Item(8, 4).Value = 0.684072732925415
$ws.Cells.Item(9, 4).Value = 0.473397970199585
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 0.5302512645721436
$ws.Cells.Item(11, 4).Value = 0.2255795747041702
$ws.Cells.Item(12, 4).Value = 0.3095614314079285
$ws.Cells.Item(13, 4).Value = 0.8774191737174988
$ws.Cells.Item(14, 4).Value = 0.8964925408363342
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 0.554570198059082
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = 0.7132629752159119
$ws.Cells.Item(17, 4).Value = 0.8830661177635193
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 0.8064461946487427
$ws.Cells.Item(19, 4).Value = 0.5176265239715576
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = 0.5133489966392517
$ws.Cells.Item(21, 4).Value = 0.2200243026018143
$ws.Cells.Item(22, 4).Value = 0.5478384494781494
$ws.Cells.Item(23, 4).Value = 0.7173118591308594
$ws.Cells.Item(24, 4).Value = 0.4508552849292755
$ws.Cells.Item(25, 4).Value = 0.8262330889701843
$ws.Cells.Item(26, 4).Value = 0.934478759765625
$ws.Cells.Item(27, 4).Value = 0.4743472039699554
$ws.Cells.Item(28, 4).Value = 0.6623991727828979
$ws.Cells.Item(29, 4).Value = 0.2277339100837708
$ws.Cells.Item(30, 4).Value = 0.6854764223098755
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 0.5086690783500671
$ws.Cells.Item(32, 4).Value = 0.6448900699615479
$ws.Cells.Item(33, 4).Value = 0.7915782928466797
$ws.Cells.Item(34, 4).Value = 0.9039463996887207
$ws.Cells.Item(35, 4).Value = 0.6719736456871033
$ws.Cells.Item(36, 4).Value = 0.2903981804847717
$ws.Cells.Item(37, 4).Value = 0.682597279548645
$ws.Cells.Item(38, 4).Value = 0.6052093505859375
$ws.Cells.Item(39, 4).Value = 0.8350102305412292
$ws.Cells.Item(40, 4).Value = 0.8173097372055054
$ws.Cells.Item(41, 4).Value = 0.9269909262657166
$ws.Cells.Item(42, 4).Value = 0.4373602569103241
$ws.Cells.Item(43, 4).Value = 0.5490264892578125
$ws.Cells.Item(44, 4).Value = 0.5116506814956665
$ws.Cells.Item(45, 4).Value = 0.3559271395206451
$ws.Cells.Item(46, 4).Value = 0.8489736318588257
$ws.Cells.Item(47, 4).Value = 0.6824671030044556
$ws.Cells.Item(48, 4).Value = 0.6174136996269226
$ws.Cells.Item(49, 4).Value = 0.4476843476295471
$ws.Cells.Item(50, 4).Value = 0.4695099294185638
$ws.Cells.Item(51, 4).Value = 0.5353732705116272
$ws.Cells.Item(52, 4).Value = 0.1604891270399094
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 0.4576481878757477
$ws.Cells.Item(54, 4).Value = 0.8174224495887756
$ws.Cells.Item(55, 4).Value = 0.3048042058944702
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(56, 4).Value = 0.5008374452590942
$ws.Cells.Item(57, 4).Value = 0.3472401797771454
$ws.Cells.Item(58, 4).Value = 0.7781393527984619
$ws.Cells.Item(59, 4).Value = 0.5788542032241821
$ws.Cells.Item(60, 4).Value = 0.8915619850158691
$ws.Cells.Item(61, 4).Value = 0.7705929279327393
$ws.Cells.Item(62, 4).Value = 0.6794452667236328
$ws.Cells.Item(63, 4).Value = 0.7870506048202515
$ws.Cells.Item(64, 4).Value = 0.3778361976146698
$ws.Cells.Item(65, 4).Value = 0.2814287841320038
$ws.Cells.Item(66, 4).Value = 0.8405328989028931
$ws.Cells.Item(67, 4).Value = 0.4315721094608307
$ws.Cells.Item(68, 3).Value = 1
$ws.Cells.Item(68, 4).Value = 0.5881571769714355
$ws.Cells.Item(69, 4).Value = 0.6431185007095337
$ws.Cells.Item(70, 4).Value = 0.2377204895019531
$ws.Cells.Item(71, 4).Value = 0.4428095817565918
$ws.Cells.Item(72, 4).Value = 0.4955498576164246
$ws.Cells.Item(73, 4).Value = 0.8746345043182373
$ws.Cells.Item(74, 4).Value = 0.450461357831955
$ws.Cells.Item(75, 4).Value = 0.2036319226026535
$ws.Cells.Item(76, 4).Value = 0.2354903817176819
$ws.Cells.Item(77, 3).Value = 1
$ws.Cells.Item(77, 4).Value = 0.6538888216018677
$ws.Cells.Item(78, 4).Value = 0.6643117070198059
$ws.Cells.Item(79, 4).Value = 0.2507842183113098
